$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add results for plant_cap_low (row 8)
$ws.Range("B8").Value = 4238322741.1799998
$ws.Range("C8").Value = 9929293.0508699995

# Add results for plant_cap_high (row 9)
$ws.Range("B9").Value = 10696284773.200001
$ws.Range("C9").Value = 0.0000019121349606199998
$ws.Range("C9").NumberFormat = "0.00E+00"

# Force recalculation of dependent formulas (D8:G8, D9:G9, B14, C14)
$excel.CalculateFull()

# Update active cell selection to C8
$ws.Range("C8").Select()
